# Update cryptos list: apply latest scraped price/volume (and some reordered
# rows) onto the existing "Sheet1" table. Cell A (rank index), header row 1,
# and all styling are left untouched - only B:E data cells for rows 2-51 change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Price column values like "1.001" / "238.36" parse as numbers through
    # plain .Value assignment, which would silently change the cell type from
    # text to numeric (and reformat "0.05200" -> 0.052, losing the trailing
    # zero). Force a Text number format for the assignment, then restore the
    # cell's style to the workbook's default 'Normal' so no stray formatting
    # is left behind, matching the source data (plain strings, no style).
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = '29.099.57'
$ws.Range("E2").Value = '  -2.73%  '
# Row 3
$ws.Range("D3").Value = '1.847.65'
$ws.Range("E3").Value = '  -1.72%  '
# Row 4
Set-TextValue "D4" '1.001'
$ws.Range("E4").Value = '  +0.23%  '
# Row 5
Set-TextValue "D5" '0.6961'
$ws.Range("E5").Value = '  -5.33%  '
# Row 6
Set-TextValue "D6" '238.36'
$ws.Range("E6").Value = '  -1.51%  '
# Row 7
Set-TextValue "D7" '1.001'
$ws.Range("E7").Value = '  +0.23%  '
# Row 8
Set-TextValue "D8" '0.3043'
$ws.Range("E8").Value = '  -3.67%  '
# Row 9
Set-TextValue "D9" '0.07518'
$ws.Range("E9").Value = '  +5.28%  '
# Row 10
Set-TextValue "D10" '23.20'
$ws.Range("E10").Value = '  -5.62%  '
# Row 11
Set-TextValue "D11" '0.08138'
$ws.Range("E11").Value = '  -1.53%  '
# Row 12
$ws.Range("D12").Value = '1.846.32'
$ws.Range("E12").Value = '  -4.09%  '
# Row 13
Set-TextValue "D13" '0.7233'
$ws.Range("E13").Value = '  -4.32%  '
# Row 14
Set-TextValue "D14" '5.207'
$ws.Range("E14").Value = '  -2.72%  '
# Row 15
Set-TextValue "D15" '88.91'
$ws.Range("E15").Value = '  -3.87%  '
# Row 16
$ws.Range("D16").Value = '29.308.62'
$ws.Range("E16").Value = '  -2.04%  '
# Row 17
Set-TextValue "D17" '5.757'
$ws.Range("E17").Value = '  -6.05%  '
# Row 18
$ws.Range("B18").Value = 'Avalanche'
$ws.Range("C18").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue "D18" '13.08'
$ws.Range("E18").Value = '  -3.04%  '
# Row 19
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue "D19" '236.31'
$ws.Range("E19").Value = '  -5.40%  '
# Row 20
Set-TextValue "D20" '0.000007667'
$ws.Range("E20").Value = '  -2.36%  '
# Row 21
$ws.Range("E21").Value = '  +0.37%  '
# Row 22
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = '2.118.43'
$ws.Range("E22").Value = '  -0.08%  '
# Row 23
$ws.Range("B23").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue "D23" '1.001'
$ws.Range("E23").Value = '  +0.32%  '
# Row 24
$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue "D24" '7.573'
$ws.Range("E24").Value = '  -2.70%  '
# Row 25
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue "D25" '9.008'
$ws.Range("E25").Value = '  -2.25%  '
# Row 26
Set-TextValue "D26" '160.03'
$ws.Range("E26").Value = '  -1.93%  '
# Row 27
$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue "D27" '0.1434'
$ws.Range("E27").Value = '  -7.10%  '
# Row 28
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue "D28" '18.06'
$ws.Range("E28").Value = '  -3.04%  '
# Row 29
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue "D29" '1.972'
$ws.Range("E29").Value = '  -3.38%  '
# Row 30
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue "D30" '1.392'
$ws.Range("E30").Value = '  -4.35%  '
# Row 31
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D31" '4.495'
$ws.Range("E31").Value = '  -1.22%  '
# Row 32
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue "D32" '1.490'
$ws.Range("E32").Value = '  -2.53%  '
# Row 33
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue "D33" '3.983'
$ws.Range("E33").Value = '  -5.13%  '
# Row 34
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D34" '0.05200'
$ws.Range("E34").Value = '  -1.96%  '
# Row 35
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue "D35" '1.185'
$ws.Range("E35").Value = '  -4.65%  '
# Row 36
$ws.Range("B36").Value = 'Frax'
$ws.Range("C36").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue "D36" '1.033'
$ws.Range("E36").Value = '  +3.64%  '
# Row 37
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D37" '0.7002'
$ws.Range("E37").Value = '  -8.14%  '
# Row 38
$ws.Range("B38").Value = 'HuobiToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue "D38" '2.665'
$ws.Range("E38").Value = '  -1.78%  '
# Row 39
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D39" '0.01858'
$ws.Range("E39").Value = '  -4.58%  '
# Row 40
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue "D40" '2.685'
$ws.Range("E40").Value = '  -2.42%  '
# Row 41
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue "D41" '0.9376'
$ws.Range("E41").Value = '  +7.40%  '
# Row 42
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue "D42" '5.987'
$ws.Range("E42").Value = '  -0.56%  '
# Row 43
$ws.Range("D43").Value = '1.068.31'
$ws.Range("E43").Value = '  -1.20%  '
# Row 44
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue "D44" '0.4273'
$ws.Range("E44").Value = '  -5.78%  '
# Row 45
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue "D45" '70.29'
$ws.Range("E45").Value = '  -2.30%  '
# Row 46
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue "D46" '1.001'
$ws.Range("E46").Value = '  +0.15%  '
# Row 47
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue "D47" '103.19'
$ws.Range("E47").Value = '  -1.64%  '
# Row 48
Set-TextValue "D48" '1.765'
$ws.Range("E48").Value = '  -4.19%  '
# Row 49
$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '2.009.92'
$ws.Range("E49").Value = '  -0.36%  '
# Row 50
$ws.Range("B50").Value = 'Aptos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue "D50" '7.029'
$ws.Range("E50").Value = '  -6.60%  '
# Row 51
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D51" '9.081'
$ws.Range("E51").Value = '  -4.50%  '
